$wb = $excel.ActiveWorkbook

# New window position/size on screen (best-effort; cosmetic window geometry)
$win = $wb.Windows.Item(1)
$win.Left = 27900
$win.Top = 1460

# Sheet1: select/activate it and move the selection to F8
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("F8").Select()

# Sheet2: new Zeta/DLS data point for C6, then switch back to it and select C7
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("C6").Value = 7000
$ws2.Range("C7").Select()
